# NOTE: worksheet object references become stale (bound to a
# positional index) once the sheet collection is mutated (Add/Move).
# To stay safe, a fresh reference is looked up by name right before
# each use, immediately after any such mutation.

$wb = $excel.ActiveWorkbook

# --- Update "Modulos" sheet (existing sheet4.xml) ---
$modulos = $wb.Worksheets.Item("Modulos")
$modulos.Range("D14").Value = "PLAZO FIJO"
$modulos.Range("D15").Value = "ARCHIVOS PERSONAS"
$modulos.Range("D16").Value = "Consulta Pagos Plazo Fijo"

# --- Add new "DataUser" sheet (headers first, to match shared-string order) ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "DataUser"

$newSheet = $wb.Worksheets.Item("DataUser")
$newSheet.Range("A1").Value = "Numero de Cuenta"
$newSheet.Range("B1").Value = "IdPersona"
$newSheet.Range("C1").Value = "NumeroOperacion"

# --- Back to "Modulos" for the last new row (re-fetch: sheet collection changed) ---
$modulos = $wb.Worksheets.Item("Modulos")
$modulos.Range("D17").Value = "Consulta de Posicion en Linea"
$modulos.Range("D17").Select()

# --- Finish "DataUser" sheet content/formatting (re-fetch again) ---
$newSheet = $wb.Worksheets.Item("DataUser")
$newSheet.Range("C2").Value = 11918739

$newSheet.Columns.Item(1).ColumnWidth = 25
$newSheet.Columns.Item(2).ColumnWidth = 14.28515625
$newSheet.Columns.Item(3).ColumnWidth = 22.42578125

$newSheet.Range("C2").Select()

# Move the new sheet to the very end, after all existing sheets
# (re-fetch once more right before the move, then stop touching it)
$newSheet = $wb.Worksheets.Item("DataUser")
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
